$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 14), matching the unformatted style of row 13's A cell
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A14").Value = 400000042

$ws.Range("BG14:CQ14").Value = 1000

# Update the view state (top-left visible cell and active selection)
$ws.Application.ActiveWindow.ScrollColumn = 81   # column CC
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("CQ13").Select()
